$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original Text format so values like
# '1.00' or '534.95' are not auto-converted into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values from the crypto data refresh
$ws.Range("D2").Value = '58.012.17'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '3.131.35'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '534.95'
$ws.Range("E5").Value = '  +2.78%  '
$ws.Range("D6").Value = '138.84'
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.504'
$ws.Range("E8").Value = '  +10.39%  '
$ws.Range("D9").Value = '7.36'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("D11").Value = '0.417'
$ws.Range("E11").Value = '  +4.23%  '
$ws.Range("E12").Value = '  +3.09%  '
$ws.Range("D13").Value = '3.670.62'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '25.66'
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("D15").Value = '0.0000170'
$ws.Range("E15").Value = '  +5.39%  '
$ws.Range("D16").Value = '58.084.64'
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("E17").Value = '  +5.85%  '
$ws.Range("D18").Value = '3.135.45'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").Value = '12.93'
$ws.Range("E19").Value = '  +3.77%  '
$ws.Range("D20").Value = '8.15'
$ws.Range("E20").Value = '  +3.70%  '
$ws.Range("D21").Value = '375.42'
$ws.Range("E21").Value = '  +7.68%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").Value = '70.09'
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("E25").Value = '  +2.84%  '
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '0.0₃0887'
$ws.Range("E28").Value = '  +2.63%  '
$ws.Range("D29").Value = '7.78'
$ws.Range("E29").Value = '  +6.83%  '
$ws.Range("E30").Value = '  +5.90%  '
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("E32").Value = '  +4.18%  '
$ws.Range("D33").Value = '5.15'
$ws.Range("E33").Value = '  +6.36%  '
$ws.Range("E34").Value = '  +4.12%  '
$ws.Range("D35").Value = '160.88'
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("D36").Value = '6.22'
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").Value = '1.34'
$ws.Range("E37").Value = '  +8.84%  '
$ws.Range("D38").Value = '25.42'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  +6.96%  '
$ws.Range("E40").Value = '  +3.01%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.595.89'
$ws.Range("E41").Value = '  +8.84%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '4.21'
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("D43").Value = '38.77'
$ws.Range("E43").Value = '  +5.81%  '
$ws.Range("D44").Value = '0.699'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("D45").Value = '0.0271'
$ws.Range("E45").Value = '  +2.74%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +4.46%  '
$ws.Range("D48").Value = '0.977'
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").Value = '0.0981'
$ws.Range("E49").Value = '  +9.52%  '
$ws.Range("D50").Value = '20.15'
$ws.Range("E50").Value = '  +2.57%  '
$ws.Range("D51").Value = '0.748'
$ws.Range("E51").Value = '  -1.98%  '
